# feat: add 2022-Q4 data
#
# Target layout after the edit:
#   Sheet1 "总计"     (unchanged name/position) - summary table gains a
#                      2022-Q4 row, the old 2020-Q4 row is pushed down.
#   Sheet2 "2022-Q4"  (NEW - inserted right after "总计") - per-fund detail
#                      for the new quarter.
#   Sheet3 "2020-Q4"  (the original Sheet2, shifted right one slot, its
#                      contents are unchanged).

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)     # "总计"

# --- 1. Insert the new "2022-Q4" worksheet right after "总计" -------------
$newQ4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $summary)
$newQ4.Name = "2022-Q4"

# --- 2. Update the "总计" summary sheet -----------------------------------
# Push the existing 2020-Q4 summary row (row 2) down to row 3, preserving
# formatting, then overwrite row 2 with the new 2022-Q4 totals.
$summary.Range("A2:D2").Copy($summary.Range("A3:D3"))

$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 0.19

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2020-Q4"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.89

# --- 3. Populate the new "2022-Q4" fund-holdings sheet ---------------------
# Copy header/cell formatting from the "总计" sheet (style used there matches
# the style the new sheet's cells should carry).
$summary.Range("B1:D1").Copy($newQ4.Range("B1:D1"))
$summary.Range("B1:D1").Copy($newQ4.Range("E1:G1"))
$summary.Range("B1").Copy($newQ4.Range("H1"))
$summary.Range("A2").Copy($newQ4.Range("A2"))
$summary.Range("A2").Copy($newQ4.Range("A3"))

$newQ4.Range("B1").Value = "基金代码"
$newQ4.Range("C1").Value = "基金名称"
$newQ4.Range("D1").Value = "基金规模"
$newQ4.Range("E1").Value = "股票总仓位"
$newQ4.Range("F1").Value = "仓位占比"
$newQ4.Range("G1").Value = "持有市值(亿元)"
$newQ4.Range("H1").Value = "仓位排名"

$newQ4.Range("A2").Value = 0
$newQ4.Range("H2").Value = 6

$newQ4.Range("A3").Value = 1
$newQ4.Range("H3").Value = 6

# The remaining data cells on rows 2/3 are stored as text in the source
# workbook (e.g. fund code "001735" must keep its leading zero), so force
# a text number-format before assigning the values, then strip the
# resulting explicit style again so the cells stay visually unstyled.
$textCells = $newQ4.Range("B2:G3")
$textCells.NumberFormat = "@"

$newQ4.Range("B2").Value = "001735"
$newQ4.Range("C2").Value = "广发百发大数据策略成长灵活配置混合E"
$newQ4.Range("D2").Value = "8.88"
$newQ4.Range("E2").Value = "90.16"
$newQ4.Range("F2").Value = "1.41"
$newQ4.Range("G2").Value = "0.1252"

$newQ4.Range("B3").Value = "001734"
$newQ4.Range("C3").Value = "广发百发大数据策略成长灵活配置混合A"
$newQ4.Range("D3").Value = "4.29"
$newQ4.Range("E3").Value = "90.16"
$newQ4.Range("F3").Value = "1.41"
$newQ4.Range("G3").Value = "0.0605"

$textCells.ClearFormats()

# --- 4. Restore the originally-active "2020-Q4" tab as the selected tab ---
$wb.Worksheets.Item(3).Select()
